# Apply updated crypto price/volume figures scraped on 2023-01-12 (GitHub Actions run).
# Source workbook stores these as literal text in columns D (Price) and E (Volume(1h));
# a leading apostrophe keeps Excel from reinterpreting the numeric-looking text as a
# number/percentage (which would silently change both the stored value and formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'284.79"
$ws.Range("E2").Value = "'2.80%"

# Row 3: OKB
$ws.Range("D3").Value = "'28.75"
$ws.Range("E3").Value = "'5.97%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'4.909"
$ws.Range("E4").Value = "'1.23%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.06487"
$ws.Range("E5").Value = "'1.29%"

# Row 6: KuCoinToken
$ws.Range("D6").Value = "'7.211"
$ws.Range("E6").Value = "'3.93%"

# Row 7: FTXToken
$ws.Range("D7").Value = "'1.338"
$ws.Range("E7").Value = "'11.50%"

# Row 8: MXToken
$ws.Range("D8").Value = "'0.9139"
$ws.Range("E8").Value = "'4.34%"

# Row 9: WazirX
$ws.Range("E9").Value = "'1.18%"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.06396"
$ws.Range("E10").Value = "'24.97%"

# Row 11: MandalaExchangeToken
$ws.Range("D11").Value = "'0.07565"
$ws.Range("E11").Value = "'0.94%"

# Row 12: BitrueCoin
$ws.Range("D12").Value = "'0.02981"
$ws.Range("E12").Value = "'0.54%"

# Row 13: BitMartToken
$ws.Range("D13").Value = "'0.08971"
$ws.Range("E13").Value = "'-0.10%"

# Row 14: BitForexToken
$ws.Range("D14").Value = "'0.001606"
$ws.Range("E14").Value = "'2.80%"

# Row 15: One
$ws.Range("D15").Value = "'0.0006561"
$ws.Range("E15").Value = "'3.43%"

# Row 16: TigerCash
$ws.Range("D16").Value = "'0.006050"
$ws.Range("E16").Value = "'-0.62%"

# Row 17: LEO
$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'-0.49%"

# Row 18: GateToken
$ws.Range("D18").Value = "'3.376"
$ws.Range("E18").Value = "'2.07%"

# Row 20: BitpandaEcosystemToken
$ws.Range("D20").Value = "'0.3148"
$ws.Range("E20").Value = "'0.40%"

# Row 21: ProBitToken
$ws.Range("D21").Value = "'0.1342"
$ws.Range("E21").Value = "'-0.43%"

# Row 22: MCDex
$ws.Range("D22").Value = "'4.012"
$ws.Range("E22").Value = "'2.50%"

# Row 23: ZBToken
$ws.Range("D23").Value = "'0.1555"
$ws.Range("E23").Value = "'12.70%"

# Row 24: CoinExToken
$ws.Range("D24").Value = "'0.04476"
$ws.Range("E24").Value = "'1.16%"

# Row 25: BitKan
$ws.Range("D25").Value = "'0.001190"
$ws.Range("E25").Value = "'1.22%"

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004323"
$ws.Range("E26").Value = "'11.87%"

# Row 28: NitroEx
$ws.Range("D28").Value = "'0.0001182"
$ws.Range("E28").Value = "'-9.08%"

# Row 29: UpBots
$ws.Range("E29").Value = "'-15.65%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.04150"
$ws.Range("E40").Value = "'-0.22%"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.006728"
$ws.Range("E41").Value = "'-1.16%"

# Row 42: BKEXToken
$ws.Range("D42").Value = "'0.1231"
$ws.Range("E42").Value = "'4.97%"

# Row 43: CEJI
$ws.Range("D43").Value = "'0.002173"
$ws.Range("E43").Value = "'11.49%"

# Row 44: LocalTraders
$ws.Range("D44").Value = "'0.01179"
$ws.Range("E44").Value = "'-1.15%"

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00005374"
$ws.Range("E45").Value = "'1.60%"

# Row 46: BOLO
$ws.Range("D46").Value = "'1.819"
$ws.Range("E46").Value = "'7.83%"

